# add support for cisco devices
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row for a Cisco Nexus device (fill the "model" cell first).
$ws.Range("A8").Value = "10.9.106.44"
$ws.Range("C8").Value = "PTSWCORE2"
$ws.Range("D8").Value = "nexus"
$ws.Range("E8").Value = "access"
$ws.Range("F8").Value = "Cisco Nexus 3048TP1GESys"
$ws.Range("G8").Value = "Node status is Up."
$ws.Range("H8").Value = "Cisco"
$ws.Range("I8").Value = "PT"

# Rename the "Type" column header to "group" and repurpose its values.
$ws.Range("B1").Value = "group"
$ws.Range("B3").Value = "none"
$ws.Range("B4").Value = "none"
$ws.Range("B5").Value = "none"
$ws.Range("B6").Value = "none"
$ws.Range("B7").Value = "none"
$ws.Range("B8").Value = "cisco"
$ws.Range("B2").Value = "junos"

# Column B now holds longer words ("group"/"backbone"-ish) - refresh its best-fit width.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Leave the selection where the author ended up after entering the new row.
$ws.Range("C15").Select() | Out-Null
